$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 128: 剑指 Offer 03. 数组中重复的数字 (数组 / 简单)
$ws.Range("A128").Value = 129
$ws.Hyperlinks.Add($ws.Range("B128"), "https://leetcode.cn/problems/shu-zu-zhong-zhong-fu-de-shu-zi-lcof/", "", "", "剑指 Offer 03. 数组中重复的数字")
$ws.Range("B128").HorizontalAlignment = -4131
$ws.Range("C128").Value = "数组"
$ws.Range("D128").Value = "简单"

# Row 129: 剑指 Offer 07. 重建二叉树 (二叉树 / 中)
$ws.Range("A129").Value = 130
$ws.Hyperlinks.Add($ws.Range("B129"), "https://leetcode.cn/problems/zhong-jian-er-cha-shu-lcof/", "", "", "剑指 Offer 07. 重建二叉树")
$ws.Range("B129").HorizontalAlignment = -4131
$ws.Range("C129").Value = "二叉树"
$ws.Range("D129").Value = "中"

$ws.Range("A1:E129").Select()
$ws.Range("B57").Select()
